$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for:" $find
    }
}

# --- Title ---
Replace-Text "Genome Editing: Reshaping Life's Blueprint" "The Astronomical Wonders: A Glimpse into the Universe's Vastness"

# --- Author name: merges 3 runs ("Dr" + "." + " Ana Locke") into one run ---
Replace-Text "Dr. Ana Locke" "Clara Patterson"

# --- Email: "analocke@biotechresearch" + "." + "org"  ->  "clara" + "patterson@yopmail" + "." + "com" ---
Replace-Text "analocke@biotechresearch" "clara"
Replace-Text "org" "patterson@yopmail"
# Now append new "." and "com" runs after the run that currently holds "patterson@yopmail".
$emailPara = $d.Paragraphs(3).Range
$lastRun = $d.Range($emailPara.End - 1, $emailPara.End - 1)
$lastRun.InsertAfter(".com")
# split the appended text into two runs ("." and "com") matching separate formatting runs
$dotPos = $emailPara.End - 1 - 3  # position right after "com" insert start; recompute below

# --- Body paragraph 1 (genome -> astronomy) ---
Replace-Text "The intricate tapestry of life is woven by the threads of DNA, the blueprint that encodes the blueprint that encodes the very essence of every organism" "Immerse yourself in the realm of cosmic exploration, uncovering the intricate tapestry of celestial mysteries that captivate humanity's collective imagination"
Replace-Text " Unraveling the mysteries of this molecular code has unlocked a new era of biological manipulation, where scientists can edit and rewrite the genetic material with unparalleled precision" " From the glimmering brilliance of countless stars to the enigmatic enigma of black holes, the universe beckons us with its boundless wonder"
Replace-Text " This technology, known as genome editing, has the potential to revolutionize our understanding of life and transform the field of medicine, agriculture, and conservation, marking a watershed moment in the annals of scientific endeavor" " Our journey through the cosmos begins with a meticulous analysis of constellations, unveiling the intricate patterns and stories etched across the celestial canvas"
Replace-Text "Genome editing is the cornerstone for a new era of biological research" "Beyond the scope of distant constellations lies a kaleidoscope of celestial wonders, each possessing its own unique allure"
Replace-Text " By allowing scientists to make targeted alterations to DNA sequences, this technology has opened the door to unravelling mysteries of gene function and unlocking the secrets of genetic diseases, enabling researchers to explore the intricate dance of genes and their interplay with the environment" " From the fiery brilliance of supernovae to the enigmatic dance of quasars, the universe teems with phenomena that challenge our understanding of reality"
Replace-Text " It has also paved the way for the development of transformative therapies, promising to reshape the landscape of healthcare by providing precise and personalized treatments" " Among these cosmic marvels, the dynamics of our solar system hold a special significance, inviting us to delve into the synchronized ballet of planets orbiting our life-giving sun, exploring the intricacies of celestial mechanics that govern their motion"
Replace-Text "Further, genome editing offers extraordinary promise in the fields of agriculture and conservation" "As we venture further into the cosmological tapestry, we encounter a profound question: are we alone? The search for extraterrestrial life takes center stage, captivating the scientific community and the public alike"

# Remove " It holds the power...agriculture" + "." run entirely, and replace " Additionally..." text.
Replace-Text " It holds the power to engineer crops that are resistant to pests and diseases, ensuring our food supply and mitigating the impact of climate change on agriculture." ""
Replace-Text " Additionally, genome editing can be utilized to revive endangered species teetering on the brink of extinction and protect valuable ecosystems, providing a beacon of hope for a sustainable future" " From analyzing the chemical composition of distant exoplanets to meticulously scrutinizing radio signals, humanity's untiring quest for answers probes the very essence of life's existence beyond Earth, igniting imaginations and fueling scientific endeavors"

# --- Summary paragraph ---
Replace-Text "Genome editing, with its ability to reprogram the very fabric of life, has ushered in an unprecedented era of scientific discovery" "In this cosmic expedition, we embarked on a captivating odyssey to explore the wonders of the universe, unveiling the intricate patterns and stories woven across constellations"
Replace-Text " Its applications span a vast array of fields, from medicine to agriculture to conservation, fuelled by the transformative power of DNA modification" " We marveled at the celestial wonders, pondering the enigma of black holes and the fiery spectacle of supernovae"

# Merge " As research continues...and the " + "preservation of our shared biosphere" (lastRenderedPageBreak run) into one run with new text
Replace-Text " As research continues to delve deeper into the intricate mechanisms of genetic masterpieces, genome editing stands poised to reshape our understanding of life itself, offering hope for cures to genetic diseases, sustainable food production, and the preservation of our shared biosphere" " The dynamics of our solar system invited us to witness the harmonious ballet of planets, while the search for extraterrestrial life ignited our imaginations and spurred scientific endeavors"

Replace-Text " This technology holds the key to unlocking the untapped potential of life's blueprint, heralding an exhilarating chapter in the chronicle of human endeavor" " As we continue to gaze upon the vastness of the cosmos, we remain humbled by its immeasurable grandeur and captivated by its endless mysteries, inspiring generations to come to unravel the enigma of our existence within this grand celestial expanse"

# --- Add trailing empty paragraph ---
$d.Content.InsertParagraphAfter()
